$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    "B2" = 4.18974298904908
    "C2" = 1.010846160492292
    "D2" = 0.3019806785874479
    "E2" = 0.04244905478493433
    "G2" = 0.002794105262723612
    "I2" = 7.428878210540375
    "J2" = 0.02118217562691882
    "L2" = 0.8093344390344441
    "M2" = 1.00447362285572
    "B3" = 4.19465038627834
    "C3" = 0.9752333885112989
    "D3" = 0.302344670193321
    "E3" = 0.04211969215485656
    "G3" = 0.002802507571648123
    "I3" = 7.125078497093824
    "J3" = 0.01844681198096509
    "L3" = 0.8049526275201515
    "M3" = 1.003509925284511
    "B4" = 4.201484049620547
    "C4" = 0.9540371283532068
    "D4" = 0.3026694991301184
    "E4" = 0.04191215212604238
    "G4" = 0.002807925473726371
    "I4" = 6.93853980787739
    "J4" = 0.01677422461466449
    "L4" = 0.8026065926257218
    "M4" = 1.003657338181924
    "B5" = 4.205226426963151
    "C5" = 0.9455662023648301
    "D5" = 0.3028273463906999
    "E5" = 0.04182622436809957
    "G5" = 0.002810198681215834
    "I5" = 6.862509587364627
    "J5" = 0.01609423538869947
    "L5" = 0.8017369195182624
    "M5" = 1.00390276037745
    "B6" = 4.205905594472995
    "C6" = 0.9441696297137696
    "D6" = 0.302855095436783
    "E6" = 0.04181187385100982
    "G6" = 0.002810580101820899
    "I6" = 6.84988359130142
    "J6" = 0.01598141648720741
    "L6" = 0.801597719469612
    "M6" = 1.003954692746817
    "B7" = 4.201530647300672
    "C7" = 0.953922213760336
    "D7" = 0.3026715247642713
    "E7" = 0.04191099877864257
    "G7" = 0.002807955865915737
    "I7" = 6.937514512358632
    "J7" = 0.01676504773351439
    "L7" = 0.8025945145828075
    "M7" = 1.003659898207118
    "B8" = 4.190640203072292
    "C8" = 0.9984266924798817
    "D8" = 0.3020851464058367
    "E8" = 0.04233658154586539
    "G8" = 0.002796948822512474
    "I8" = 7.324120989000761
    "J8" = 0.02023749641305983
    "L8" = 0.8077519501703847
    "M8" = 1.003987617503093
    "B9" = 4.199750640960133
    "C9" = 1.091104681507034
    "D9" = 0.3017398586266395
    "E9" = 0.04312986625766779
    "G9" = 0.002777405148661663
    "I9" = 8.082848920632671
    "J9" = 0.02710873624928212
    "L9" = 0.820611548326923
    "M9" = 1.010520663586277
    "B10" = 4.225248029839179
    "C10" = 1.162621968421831
    "D10" = 0.3019778784772598
    "E10" = 0.04368874468657324
    "G10" = 0.002764272775331065
    "I10" = 8.641576544819941
    "J10" = 0.03220460023737104
    "L10" = 0.8317542251247261
    "M10" = 1.018951175908768
    "B11" = 4.240983897981494
    "C11" = 1.195930980736421
    "D11" = 0.3021932553595548
    "E11" = 0.04393806216136831
    "G11" = 0.002758560941591468
    "I11" = 8.896243519053371
    "J11" = 0.03453538431614334
    "L11" = 0.8371961568799975
    "M11" = 1.023584063744948
    "B12" = 4.24754193731394
    "C12" = 1.20865808798078
    "D12" = 0.3022902389923985
    "E12" = 0.04403178646069339
    "G12" = 0.002756435414589686
    "I12" = 8.992767199198056
    "J12" = 0.03541998597848561
    "L12" = 0.8393108987461062
    "M12" = 1.025453891070917
    "B13" = 4.246102825732578
    "C13" = 1.205911984880061
    "D13" = 0.3022686653854123
    "E13" = 0.04401163147138032
    "G13" = 0.00275689152508014
    "I13" = 8.971974992852552
    "J13" = 0.03522938025906797
    "L13" = 0.8388530437877932
    "M13" = 1.02504604353301
    "B14" = 4.241511397501824
    "C14" = 1.196975755825463
    "D14" = 0.3022009249717428
    "E14" = 0.04394578656164505
    "G14" = 0.002758385324822711
    "I14" = 8.904182752802399
    "J14" = 0.03460812029808835
    "L14" = 0.837369053642476
    "M14" = 1.023735577421235
    "B15" = 4.238777180579916
    "C15" = 1.191516931877402
    "D15" = 0.3021614415757483
    "E15" = 0.0439053658256725
    "G15" = 0.002759305185393648
    "I15" = 8.862669807913676
    "J15" = 0.03422784402012979
    "L15" = 0.836467110006339
    "M15" = 1.022947936688787
    "B16" = 4.224303265461799
    "C16" = 1.160460936622371
    "D16" = 0.3019659602340568
    "E16" = 0.04367235396581837
    "G16" = 0.002764651307977676
    "I16" = 8.624944781251145
    "J16" = 0.03205254660902312
    "L16" = 0.8314061191049831
    "M16" = 1.018664521368649
    "B17" = 4.216486583966741
    "C17" = 1.14160924417655
    "D17" = 0.3018734861441317
    "E17" = 0.04352816333378229
    "G17" = 0.002767997927616688
    "I17" = 8.479246605175831
    "J17" = 0.03072142848303372
    "L17" = 0.828397164509866
    "M17" = 1.016241631578936
    "B18" = 4.212379759330645
    "C18" = 1.130839102104687
    "D18" = 0.3018303766952783
    "E18" = 0.04344476443294454
    "G18" = 0.002769947505687516
    "I18" = 8.395490914187434
    "J18" = 0.02995698511191591
    "L18" = 0.8267015867396026
    "M18" = 1.014923097210321
    "B19" = 4.211055962631633
    "C19" = 1.127204963567749
    "D19" = 0.3018175109962016
    "E19" = 0.04341644664103494
    "G19" = 0.00277061184848548
    "I19" = 8.367140130568828
    "J19" = 0.02969835551167677
    "L19" = 0.8261335100302176
    "M19" = 1.014489531898427
    "B20" = 4.217278379510276
    "C20" = 1.143608483290564
    "D20" = 0.3018822868435223
    "E20" = 0.04354356057449582
    "G20" = 0.002767639120629332
    "I20" = 8.494751550914316
    "J20" = 0.03086300491060001
    "L20" = 0.8287138383261237
    "M20" = 1.016491779988627
    "B21" = 4.242843713712659
    "C21" = 1.199597439099534
    "D21" = 0.3022204031362037
    "E21" = 0.04396514528214901
    "G21" = 0.002757945546381185
    "I21" = 8.924092498994526
    "J21" = 0.03479054421871552
    "L21" = 0.8378034690785654
    "M21" = 1.024117353746973
    "B22" = 4.26304690167143
    "C22" = 1.236853296082586
    "D22" = 0.3025313039223931
    "E22" = 0.04423668177048157
    "G22" = 0.002751828219757289
    "I22" = 9.205204920225412
    "J22" = 0.03736906373196547
    "L22" = 0.8440589983345745
    "M22" = 1.029774417580256
    "B23" = 4.251942838825073
    "C23" = 1.21690765987006
    "D23" = 0.3023571337712667
    "E23" = 0.04409211615364583
    "G23" = 0.002755073297874174
    "I23" = 9.055117903958887
    "J23" = 0.03599173818153645
    "L23" = 0.8406913680314574
    "M23" = 1.02669327902727
    "B24" = 4.216919203371958
    "C24" = 1.142704415813284
    "D24" = 0.3018782767301218
    "E24" = 0.04353660104648682
    "G24" = 0.002767801257616443
    "I24" = 8.487741740754927
    "J24" = 0.03079899562229826
    "L24" = 0.8285705632512617
    "M24" = 1.0163784561543
    "B25" = 4.194001074140544
    "C25" = 1.065440776668368
    "D25" = 0.3017470150730261
    "E25" = 0.04291955643820611
    "G25" = 0.002782475569276754
    "I25" = 7.877434371012697
    "J25" = 0.02524225850285688
    "L25" = 0.8093344390344441
    "M25" = 1.008118686976751
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value2 = $values[$cell]
}
